# Fruta / hortaliza, semanal
# Insert a new weekly record as row 60 (pushing the existing rows 60-71
# down to 61-72), populated with the week's price data for
# Pomelo / Start Ruby at Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 60:71 down to 61:72, leaving a blank row 60 to fill in.
$ws.Rows(60).Insert()

$ws.Cells.Item(60, 1).Value  = 9
$ws.Cells.Item(60, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(60, 3).Value  = "Metropolitana"
$ws.Cells.Item(60, 4).Value  = 44722
$ws.Cells.Item(60, 5).Value  = 13
$ws.Cells.Item(60, 6).Value  = "Fruta"
$ws.Cells.Item(60, 7).Value  = 100102
$ws.Cells.Item(60, 8).Value  = "Cítricos"
$ws.Cells.Item(60, 9).Value  = 100102006
$ws.Cells.Item(60, 10).Value = "Pomelo"
$ws.Cells.Item(60, 11).Value = "Start Ruby"
$ws.Cells.Item(60, 12).Value = "Primera"
$ws.Cells.Item(60, 13).Value = 380
$ws.Cells.Item(60, 14).Value = 7500
$ws.Cells.Item(60, 15).Value = 7500
$ws.Cells.Item(60, 16).Value = 7500
$ws.Cells.Item(60, 17).Value = "$/caja 14 kilos"
$ws.Cells.Item(60, 18).Value = "Región Metropolitana"
$ws.Cells.Item(60, 19).Value = 536
$ws.Cells.Item(60, 20).Value = 14
